# Regen save_data: recompute the K column (strike count, formerly "Strike#")
# from the underlying std/mean calc and write the resulting s_vals back into
# column G for every data row (rows 2-73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# s_vals: newly calculated K values for rows 2..73 (row 54 unchanged at 0)
$sVals = @(2,2,2,1,0,0,1,1,1,1,1,1,2,0,2,1,2,1,0,0,1,2,0,1,0,3,3,1,0,1,1,0,0,2,1,1,1,1,1,2,0,1,1,1,2,1,2,2,1,1,1,1,0,2,1,1,0,2,1,0,2,1,1,1,2,1,0,0,2,0,0,2)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
